$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.076720952987671
$ws.Range("B1").Value = 4.632126331329346
$ws.Range("C1").Value = 3.729093313217163
$ws.Range("D1").Value = 3.265719413757324
$ws.Range("E1").Value = 1.334676504135132
